$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: quantity 1 -> 2
$ws.Range("D5").Value = "2"

# Row 6: product 1005153 -> 1021370
$ws.Range("A6").Value = "1021370"
$ws.Range("B6").Value = "Glina. JET Cookies&Cream 18Plegx6undx50g"

# Row 7: product 1005147 -> 1033837
$ws.Range("A7").Value = "1033837"
$ws.Range("B7").Value = "Glina. JET Cookies&Cream 6Plegx18unx50g"

# Row 8: product 1001485 -> 1034056, quantity 2 -> 1
$ws.Range("A8").Value = "1034056"
$ws.Range("B8").Value = "Glina. JET Cookies&Cream EXH60Bjx2unx50g"
$ws.Range("D8").Value = "1"
